# Add new test rows (47-80) to the worksheet, covering additional data
# source/data-type combinations (Oracle, SQL Server, ClickHouse) alongside
# the existing cases.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@("snc-connector-test-oracle-string-1", "good request, data retrieved (no scheme check)", "id='70000'", "Work_Center", 0, 0, 0)
    ,@("snc-connector-test-oracle-string-2", "good request, data retrieved (no scheme check)", "product_type='PCBL'", "Product", 0, 0, 0)
    ,@("snc-connector-test-oracle-string-3", "good request, data retrieved (no scheme check)", "sales_order='UPM'", "Product_Order", 0, 0, 0)
    ,@("snc-connector-test-oracle-string-4", "good request, data retrieved (no scheme check)", "product_level='ME'", "Production_Procedure", 0, 0, 0)
    ,@("snc-connector-test-oracle-string-5", "good request, data retrieved (no scheme check)", "procedure_id='A0001'", "Work_Position", 0, 0, 0)
    ,@("snc-connector-test-oracle-string-6", "good request, data retrieved (no scheme check)", "val='45'", "KPI", 0, 0, 0)
    ,@("snc-connector-test-oracle-string-7", "good request, data retrieved (no scheme check)", "registered_capital='50000000'", "Plant_Owner", 0, 0, 0)
    ,@("snc-connector-test-oracle-string-8", "good request, data retrieved (no scheme check)", "test_result='P'", "Product_Qualification", 0, 0, 0)
    ,@("snc-connector-test-oracle-number-1", "good request, data retrieved (no scheme check)", "id='3'", "Device", 0, 0, 0)
    ,@("snc-connector-test-oracle-float-1", "good request, data retrieved (no scheme check)", "amps='23.3'", "Device", 0, 0, 0)
    ,@("snc-connector-test-oracle-time-1", "good request, data retrieved (no scheme check)", "start_time>'2019-08-07 00:00:00'", "Preactor_Order", 0, 0, 0)
    ,@("snc-connector-test-oracle-time-2", "good request, data retrieved (no scheme check)", "end_time<'2019-08-08 00:00:00'", "Preactor_Order", 0, 0, 0)
    ,@("snc-connector-test-oracle-time-3", "good request, data retrieved (no scheme check)", "start_time>'2019-08-07T00:00:00'", "Preactor_Order", 0, 0, 0)
    ,@("snc-connector-test-oracle-time-4", "good request, data retrieved (no scheme check)", "end_time<'2019-08-08T00:00:00'", "Preactor_Order", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-string-1", "good request, data retrieved (no scheme check)", "product='A5E33593642'", "Product_Order", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-float-1", "good request, data retrieved (no scheme check)", "quantity='25.0'", "Product_Order", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-bit-1", "good request, data retrieved (no scheme check)", "IgnoreShortages='false'", "BOM", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-bit-2", "good request, data retrieved (no scheme check)", "IgnoreShortages='true'", "BOM", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-bit-3", "good request, data retrieved (no scheme check)", "IgnoreShortages=false", "BOM", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-bit-4", "good request, data retrieved (no scheme check)", "IgnoreShortages=true", "BOM", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-bit-5", "good request, data retrieved (no scheme check, no condition check)", "IgnoreShortages=0", "BOM", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-bit-6", "good request, data retrieved (no scheme check, no condition check)", "IgnoreShortages=1", "BOM", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-int-1", "good request, data retrieved (no scheme check)", "consume_quantity='1'", "BOM", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-time-1", "good request, data retrieved (no scheme check)", "demand_date>'2019-08-13 07:23:50'", "Product_Order", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-time-2", "good request, data retrieved (no scheme check)", "demand_date<'2019-08-13 07:23:50'", "Product_Order", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-time-3", "good request, data retrieved (no scheme check)", "demand_date>'2019-08-13T07:23:50'", "Product_Order", 0, 0, 0)
    ,@("snc-connector-test-sqlserver-time-4", "good request, data retrieved (no scheme check)", "demand_date<'2019-08-13T07:23:50'", "Product_Order", 0, 0, 0)
    ,@("snc-connector-test-clickhouse-string-1", "good request, data retrieved (no scheme check)", "id='SIMANTIC300'", "SINAMICS_300_Log", 1, 2, 0)
    ,@("snc-connector-test-clickhouse-int-1", "good request, data retrieved (no scheme check)", "port='8'", "SINAMICS_300_Log", 1, 2, 0)
    ,@("snc-connector-test-clickhouse-float-1", "good request, data retrieved (no scheme check)", "outputcurrent_actual_AI0='98.21999'", "SINAMICS_300_Log", 0, 0, 0)
    ,@("snc-connector-test-clickhouse-time-1", "good request, data retrieved (no scheme check)", "update_time>'2021-06-03 04:50:19'", "SINAMICS_300_Log", 1, 2, 0)
    ,@("snc-connector-test-clickhouse-time-2", "good request, data retrieved (no scheme check)", "update_time<'2021-06-03 04:50:19'", "SINAMICS_300_Log", 1, 2, 0)
    ,@("snc-connector-test-clickhouse-time-3", "good request, data retrieved (no scheme check)", "update_time>'2021-06-03T04:50:19'", "SINAMICS_300_Log", 1, 2, 0)
    ,@("snc-connector-test-clickhouse-time-4", "good request, data retrieved (no scheme check)", "update_time<'2021-06-03T04:50:19'", "SINAMICS_300_Log", 1, 2, 0)
)

$startRow = 47
for ($idx = 0; $idx -lt $rows.Count; $idx++) {
    $r = $startRow + $idx
    $data = $rows[$idx]

    $ws.Cells.Item($r, 1).Value = $data[0]

    # Column B ("description") reuses the bordered style already used by the
    # header/description column in the existing rows.
    $ws.Range("B2").Copy() | Out-Null
    $ws.Cells.Item($r, 2).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 2).Value = $data[1]

    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 6).Value = $data[3]
    $ws.Cells.Item($r, 8).Value = $data[4]
    $ws.Cells.Item($r, 9).Value = $data[5]
    $ws.Cells.Item($r, 12).Value = $data[6]
}

$excel.CutCopyMode = 0
